# Apply updated crypto price/volume data (scraped refresh) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.201.37'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '3.339.50'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = "'586.64"
$ws.Range('E5').Value = '  +5.33%  '
$ws.Range('D6').Value = "'185.48"
$ws.Range('E6').Value = '  -1.10%  '
$ws.Range('E8').Value = '  -1.81%  '
$ws.Range('E9').Value = '  -0.96%  '
$ws.Range('D10').Value = "'0.582"
$ws.Range('E10').Value = '  -0.88%  '
$ws.Range('D11').Value = "'47.02"
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('E12').Value = '  -0.48%  '
$ws.Range('D13').Value = "'666.63"
$ws.Range('E13').Value = '  +10.35%  '
$ws.Range('D14').Value = '3.874.87'
$ws.Range('E14').Value = '  +0.70%  '
$ws.Range('D15').Value = "'8.48"
$ws.Range('E15').Value = '  -2.71%  '
$ws.Range('D16').Value = '66.361.60'
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.343.18'
$ws.Range('E18').Value = '  +0.58%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = "'17.87"
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('D20').Value = "'11.09"
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('E21').Value = '  -1.41%  '
$ws.Range('D22').Value = "'17.68"
$ws.Range('E22').Value = '  -4.28%  '
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').Value = "'5.04"
$ws.Range('E23').Value = '  -1.49%  '
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = "'100.81"
$ws.Range('E24').Value = '  +0.75%  '
$ws.Range('D25').Value = "'4.02"
$ws.Range('E25').Value = '  +1.45%  '
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('D27').Value = "'9.46"
$ws.Range('E27').Value = '  -1.58%  '
$ws.Range('D28').Value = "'32.29"
$ws.Range('E28').Value = '  +6.40%  '
$ws.Range('D29').Value = "'8.49"
$ws.Range('E29').Value = '  -1.61%  '
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('D31').Value = "'611.02"
$ws.Range('E31').Value = '  +4.72%  '
$ws.Range('D33').Value = "'11.06"
$ws.Range('E33').Value = '  -0.63%  '
$ws.Range('D34').Value = '3.879.62'
$ws.Range('E34').Value = '  +4.59%  '
$ws.Range('E35').Value = '  +0.24%  '
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').Value = "'56.33"
$ws.Range('E37').Value = '  -1.73%  '
$ws.Range('E38').Value = '  -2.06%  '
$ws.Range('E39').Value = '  -2.98%  '
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').Value = "'32.97"
$ws.Range('E40').Value = '  -2.77%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').Value = "'2.66"
$ws.Range('E41').Value = '  -0.39%  '
$ws.Range('E42').Value = '  -2.80%  '
$ws.Range('D43').Value = "'3.41"
$ws.Range('E43').Value = '  +1.77%  '
$ws.Range('D44').Value = "'0.336"
$ws.Range('E44').Value = '  -1.61%  '
$ws.Range('E45').Value = '  -1.14%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').Value = "'0.128"
$ws.Range('E46').Value = '  -1.43%  '
$ws.Range('B47').Value = 'CoreDAO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D47').Value = "'2.94"
$ws.Range('E47').Value = '  -15.74%  '
$ws.Range('E48').Value = '  +0.22%  '
$ws.Range('D49').Value = "'2.55"
$ws.Range('E49').Value = '  -1.60%  '
$ws.Range('E50').Value = '  +3.88%  '
$ws.Range('D51').Value = "'129.03"
$ws.Range('E51').Value = '  +4.53%  '
